# Apply cryptocurrency price/volume updates per the target diff.
# D (Price) and E (Volume(1h)) columns hold number-like text (e.g. "1.00",
# "71.248.85", "  +0.27%  ") that Excel would otherwise auto-coerce into
# numbers/dates and silently mangle (e.g. "1.00" -> 1, losing the trailing
# zero formatting and the cell type). Forcing NumberFormat "@" (Text) on
# those cells before assignment keeps the literal string intact, matching
# the inline string content in the source workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.248.85'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.812.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '704.12'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.10'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.809.52'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.45'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.57'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.458.15'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.802.52'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.464.54'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.52'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '510.92'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.46'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.716'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.89'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000141'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.74'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.960.05'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.31'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.20%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.01'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.02'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.92%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.26'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.38'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.13'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.82%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.96%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.780.28'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.62'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +10.37%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.40'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.97%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.13%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.22'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '166.25'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '49.96'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.35%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '430.60'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000303'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.92%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.57'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +8.53%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.40'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.08%  '
